# Progress Report - Group A.xlsx
# Add a "Week 5" column (F) to the weekly status table (rows 16-24) and a
# matching "Week 5" note in the notes row (row 28). Also clears the stray
# empty/formatted cell that used to sit at H16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell for the weekly table, formatted like the others --
$ws.Range("F16").Value = "Week 5"
$ws.Range("E16").Copy() | Out-Null
$ws.Range("F16").PasteSpecial(-4122) | Out-Null

# --- Same task text repeated for every team member on week 5 ----------
$weekFiveTask = "Dataset analysis, Missing values identification"
$ws.Range("F17").Value = $weekFiveTask
$ws.Range("F18").Value = $weekFiveTask
$ws.Range("F19").Value = $weekFiveTask
$ws.Range("F20").Value = $weekFiveTask
$ws.Range("F21").Value = $weekFiveTask
$ws.Range("F22").Value = $weekFiveTask
$ws.Range("F23").Value = $weekFiveTask
$ws.Range("F24").Value = $weekFiveTask

# --- Size column F to fit the week-5 table contents (matches how C:E were
#     auto-fit earlier) before the longer notes text is added to row 28 --
$ws.Columns.Item(6).AutoFit() | Out-Null

# --- Notes row: describe what happened on week 5, formatted like the
#     other notes cells (B28:D28 use wrap/vertical-top, no indent) -------
$ws.Range("F28").Value = "Week 5, Dataset analysis was continued while learning how to use different python libraries. Target variable imbalance was identified on the dataset."
$ws.Range("D28").Copy() | Out-Null
$ws.Range("F28").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Drop the stray formatted-but-empty cell that used to sit at H16 --
$ws.Range("H16").Clear() | Out-Null

# --- Restore the normal (non-scrolled) view with a fresh selection ----
$ws.Range("H24").Select() | Out-Null

Write-Output "done"
